# Auto-generated edit script applying numeric corrections per the commit diff
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 3919
$ws.Range("J17").Value = 3919
$ws.Range("L17").Value = 11757
$ws.Range("N17").Value = -12093
$ws.Range("H40").Value = 2542.2
$ws.Range("J40").Value = 4002
$ws.Range("L40").Value = 4002
$ws.Range("N40").Value = -4352
$ws.Range("H41").Value = 565.2143
$ws.Range("I41").Value = 178.63637
$ws.Range("K41").Value = 178.63637
$ws.Range("M41").Value = 261.36363
$ws.Range("H99").Value = 446.57144
$ws.Range("J99").Value = 0
$ws.Range("L99").Value = 0
$ws.Range("N99").ClearContents()
$ws.Range("H100").Value = 2511.0588
$ws.Range("I100").Value = 1637
$ws.Range("J100").Value = 3288
$ws.Range("K100").Value = 1637
$ws.Range("L100").Value = 3288
$ws.Range("M100").Value = -1096
$ws.Range("N100").Value = -4370
$ws.Range("H138").Value = 3358.7026
$ws.Range("J138").Value = 3508.75
$ws.Range("L138").Value = 10526.25
$ws.Range("N138").Value = -20806.25

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 6808.2197
$ws.Range("I132").Value = 4958.871
$ws.Range("K132").Value = 14876.613
$ws.Range("M132").Value = -12346.613

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 254140.25
$ws.Range("I86").Value = 2531
$ws.Range("K86").Value = 2531
$ws.Range("M86").Value = -1408
$ws.Range("H89").Value = 254140.25
$ws.Range("I89").Value = 2531
$ws.Range("K89").Value = 12655
$ws.Range("M89").Value = -7039
$ws.Range("H134").Value = 6495.564
$ws.Range("I134").Value = 5032.52
$ws.Range("K134").Value = 15097.56
$ws.Range("M134").Value = -12562.56

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H20").Value = 77356
$ws.Range("J20").Value = 77356
$ws.Range("L20").Value = 77356
$ws.Range("N20").Value = -77828
$ws.Range("H30").Value = 77356
$ws.Range("J30").Value = 77356
$ws.Range("L30").Value = 77356
$ws.Range("N30").Value = -77538
$ws.Range("H31").Value = 13893122
$ws.Range("J31").Value = 6537.8945
$ws.Range("L31").Value = 6537.8945
$ws.Range("N31").Value = -7127.8945
$ws.Range("H33").Value = 1674
$ws.Range("I33").Value = 1065.3334
$ws.Range("J33").Value = 3500
$ws.Range("K33").Value = 1065.3334
$ws.Range("L33").Value = 3500
$ws.Range("M33").Value = -686.3334
$ws.Range("N33").Value = -4258
$ws.Range("H34").Value = 13893122
$ws.Range("J34").Value = 6537.8945
$ws.Range("L34").Value = 6537.8945
$ws.Range("N34").Value = -6941.8945
$ws.Range("H44").Value = 12017.5
$ws.Range("I44").Value = 9064
$ws.Range("J44").Value = 14971
$ws.Range("K44").Value = 9064
$ws.Range("L44").Value = 14971
$ws.Range("M44").Value = -8622
$ws.Range("N44").Value = -15855
$ws.Range("H128").Value = 77356
$ws.Range("J128").Value = 77356
$ws.Range("L128").Value = 77356
$ws.Range("N128").Value = -87316
$ws.Range("H132").Value = 69046.44500000001
$ws.Range("I132").Value = 4373.75
$ws.Range("J132").Value = 120784.6
$ws.Range("K132").Value = 13121.25
$ws.Range("L132").Value = 362353.8
$ws.Range("M132").Value = -10591.25
$ws.Range("N132").Value = -367413.8

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H6").Value = 212.28572
$ws.Range("I6").Value = 153.14285
$ws.Range("K6").Value = 459.42855
$ws.Range("M6").Value = -346.42855
$ws.Range("H16").Value = 416.66666
$ws.Range("J16").Value = 200
$ws.Range("L16").Value = 600
$ws.Range("N16").Value = -946
$ws.Range("H39").Value = 5763.7
$ws.Range("I39").Value = 1000
$ws.Range("K39").Value = 3000
$ws.Range("M39").Value = -2706
$ws.Range("H55").Value = 4324.2
$ws.Range("I55").Value = 977.6
$ws.Range("J55").Value = 5997.5
$ws.Range("K55").Value = 2932.8
$ws.Range("L55").Value = 17992.5
$ws.Range("M55").Value = -2755.8
$ws.Range("N55").Value = -18346.5
$ws.Range("H86").Value = 297.2143
$ws.Range("I86").Value = 244.375
$ws.Range("J86").Value = 367.66666
$ws.Range("K86").Value = 733.125
$ws.Range("L86").Value = 1102.99998
$ws.Range("M86").Value = 452.875
$ws.Range("N86").Value = -3474.99998
$ws.Range("H89").Value = 297.2143
$ws.Range("I89").Value = 244.375
$ws.Range("J89").Value = 367.66666
$ws.Range("K89").Value = 2199.375
$ws.Range("L89").Value = 3308.99994
$ws.Range("M89").Value = 3728.625
$ws.Range("N89").Value = -15164.99994
$ws.Range("H129").Value = 26317054
$ws.Range("I129").Value = 587
$ws.Range("J129").Value = 71431000
$ws.Range("K129").Value = 1761
$ws.Range("L129").Value = 214293000
$ws.Range("M129").Value = 3239
$ws.Range("N129").Value = -214303000
$ws.Range("H131").Value = 8550707
$ws.Range("I131").Value = 23810516
$ws.Range("K131").Value = 71431548
$ws.Range("M131").Value = -71426508

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H31").Value = 1999
$ws.Range("I31").Value = 1999
$ws.Range("K31").Value = 1999
$ws.Range("M31").Value = -1707
$ws.Range("H37").Value = 1999
$ws.Range("I37").Value = 1999
$ws.Range("K37").Value = 1999
$ws.Range("M37").Value = -1722
$ws.Range("H70").Value = 6848.4287
$ws.Range("J70").Value = 8954.5
$ws.Range("L70").Value = 8954.5
$ws.Range("N70").Value = -9494.5
$ws.Range("H73").Value = 6848.4287
$ws.Range("J73").Value = 8954.5
$ws.Range("L73").Value = 8954.5
$ws.Range("N73").Value = -10826.5
$ws.Range("H97").Value = 1143.7188
$ws.Range("I97").Value = 1025.7273
$ws.Range("J97").Value = 1403.3
$ws.Range("K97").Value = 1025.7273
$ws.Range("L97").Value = 1403.3
$ws.Range("M97").Value = -529.7273
$ws.Range("N97").Value = -2395.3
$ws.Range("H113").Value = 4751
$ws.Range("I113").Value = 3627
$ws.Range("J113").Value = 6999
$ws.Range("K113").Value = 3627
$ws.Range("L113").Value = 6999
$ws.Range("M113").Value = -1457
$ws.Range("N113").Value = -11339
$ws.Range("H132").Value = 17000
$ws.Range("I132").Value = 0
$ws.Range("J132").Value = 17000
$ws.Range("K132").Value = 0
$ws.Range("L132").Value = 51000
$ws.Range("M132").ClearContents()
$ws.Range("N132").Value = -56060

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 4214.552
$ws.Range("I22").Value = 1955.3334
$ws.Range("J22").Value = 5231.2
$ws.Range("K22").Value = 1955.3334
$ws.Range("L22").Value = 5231.2
$ws.Range("M22").Value = -1660.3334
$ws.Range("N22").Value = -5821.2
$ws.Range("H27").Value = 4214.552
$ws.Range("I27").Value = 1955.3334
$ws.Range("J27").Value = 5231.2
$ws.Range("K27").Value = 1955.3334
$ws.Range("L27").Value = 5231.2
$ws.Range("M27").Value = -1848.3334
$ws.Range("N27").Value = -5445.2
$ws.Range("H46").Value = 8995.105
$ws.Range("J46").Value = 10411.9375
$ws.Range("L46").Value = 10411.9375
$ws.Range("N46").Value = -10787.9375
$ws.Range("H132").Value = 5831.4917
$ws.Range("I132").Value = 4779.7383
$ws.Range("K132").Value = 14339.2149
$ws.Range("M132").Value = -11809.2149
$ws.Range("H136").Value = 336316.94
$ws.Range("I136").Value = 483255.6
$ws.Range("K136").Value = 1449766.8
$ws.Range("M136").Value = -1447216.8
$ws.Range("H140").Value = 101499
$ws.Range("I140").Value = 87000
$ws.Range("J140").Value = 104398.8
$ws.Range("K140").Value = 87000
$ws.Range("L140").Value = 104398.8
$ws.Range("M140").Value = -81820
$ws.Range("N140").Value = -114758.8

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 4975.25
$ws.Range("I132").Value = 4070.9395
$ws.Range("J132").Value = 6964.7334
$ws.Range("K132").Value = 12212.8185
$ws.Range("L132").Value = 20894.2002
$ws.Range("M132").Value = -9682.818499999999
$ws.Range("N132").Value = -25954.2002
$ws.Range("H136").Value = 2730.0908
$ws.Range("I136").Value = 2055.65
$ws.Range("J136").Value = 3767.6924
$ws.Range("K136").Value = 6166.950000000001
$ws.Range("L136").Value = 11303.0772
$ws.Range("M136").Value = -3616.950000000001
$ws.Range("N136").Value = -16403.0772
